# Add a new team "Arran Katoko FC" (time_id 19833277) to the ranking table.
# The table is sorted by time_id ascending, so the new team belongs right
# before the existing row for time_id 20651178 ("Pontaç0 F.C."), which is
# currently row 18. The three rows below it (Pontaç0 F.C., lsauer fc,
# Grêmio_Campeão_LA_27) shift down by one row, and a new last row (21) is
# used for the team that used to be in row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 18-20 down to 19-21, working bottom-up so we never
# overwrite a value before it has been copied onward.
$ws.Cells.Item(21, 1).Value2 = $ws.Cells.Item(20, 1).Value2
$ws.Cells.Item(21, 2).Value2 = $ws.Cells.Item(20, 2).Value2
$ws.Cells.Item(21, 3).Value2 = $ws.Cells.Item(20, 3).Value2

$ws.Cells.Item(20, 1).Value2 = $ws.Cells.Item(19, 1).Value2
$ws.Cells.Item(20, 2).Value2 = $ws.Cells.Item(19, 2).Value2
$ws.Cells.Item(20, 3).Value2 = $ws.Cells.Item(19, 3).Value2

$ws.Cells.Item(19, 1).Value2 = $ws.Cells.Item(18, 1).Value2
$ws.Cells.Item(19, 2).Value2 = $ws.Cells.Item(18, 2).Value2
$ws.Cells.Item(19, 3).Value2 = $ws.Cells.Item(18, 3).Value2

# Write the new team into the now-vacated row 18.
$ws.Cells.Item(18, 1).Value2 = 19833277
$ws.Cells.Item(18, 2).Value2 = "Arran Katoko FC"
$ws.Cells.Item(18, 3).Value2 = 0

# Row 21's time_id cell (column A) needs the same bold/centered/bordered
# style used by the rest of column A; copy it from a neighboring cell.
$ws.Cells.Item(17, 1).Copy()
$ws.Cells.Item(21, 1).PasteSpecial(-4122)
